$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 3..9 (x=5..100). We need to insert six
# new rows so that the existing row 8 (x=50) ends up at row 10, and the
# existing row 9 (x=100) ends up at row 15, and fill in the newly inserted
# rows with the additional (x, ratio) samples.

# Step 1: push old row 8 (x=50) down by two rows -> becomes row 10.
$ws.Rows("8:9").Insert()

# Step 2: push old row 9 (x=100, now sitting at row 11) down by four more
# rows -> becomes row 15.
$ws.Rows("11:14").Insert()

# New row 8: x=30
$ws.Range("A8").Value = 30
$ws.Range("C8").Value = 45.79

# New row 9: x=40
$ws.Range("A9").Value = 40
$ws.Range("C9").Value = 55.78

# Row 10 already holds the old row-8 data (x=50) after the inserts.

# New row 11: x=60
$ws.Range("A11").Value = 60
$ws.Range("C11").Value = 70.51

# New row 12: x=70
$ws.Range("A12").Value = 70
$ws.Range("C12").Value = 75.74

# New row 13: x=80
$ws.Range("A13").Value = 80
$ws.Range("C13").Value = 79.72

# New row 14: x=90
$ws.Range("A14").Value = 90
$ws.Range("C14").Value = 83.1

# Row 15 already holds the old row-9 data (x=100) after the inserts, but
# the C value changed from 85.45 to 85.47.
$ws.Range("C15").Value = 85.47
